$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells that would otherwise be auto-converted to numbers by Excel
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '30.507.73'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").Value = '1.851.40'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '233.75'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4710'
$ws.Range("E7").Value = '  -0.57%  '
$ws.Range("D8").Value = '0.2741'
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.06340'
$ws.Range("E9").Value = '  -1.43%  '
$ws.Range("D10").Value = '17.63'
$ws.Range("E10").Value = '  +7.82%  '
$ws.Range("D11").Value = '1.847.48'
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D12").Value = '0.07407'
$ws.Range("E12").Value = '  -0.78%  '
$ws.Range("D13").Value = '5.053'
$ws.Range("E13").Value = '  +1.26%  '
$ws.Range("D14").Value = '84.57'
$ws.Range("E14").Value = '  -1.22%  '
$ws.Range("D15").Value = '0.6253'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").Value = '30.473.87'
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").Value = '241.97'
$ws.Range("E17").Value = '  +5.05%  '
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("D19").Value = '12.67'
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").Value = '0.000007339'
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").Value = '4.935'
$ws.Range("E22").Value = '  -1.30%  '
$ws.Range("D23").Value = '5.972'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").Value = '9.231'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("D25").Value = '161.87'
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("D26").Value = '18.00'
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").Value = '1.883'
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("E28").Value = '  -2.54%  '
$ws.Range("D29").Value = '1.361'
$ws.Range("E29").Value = '  -2.94%  '
$ws.Range("D30").Value = '4.017'
$ws.Range("E30").Value = '  -3.24%  '
$ws.Range("D31").Value = '3.839'
$ws.Range("E31").Value = '  -2.37%  '
$ws.Range("D32").Value = '0.04866'
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("D33").Value = '1.136'
$ws.Range("E33").Value = '  -2.60%  '
$ws.Range("D34").Value = '0.7056'
$ws.Range("E34").Value = '  -2.64%  '
$ws.Range("D35").Value = '2.712'
$ws.Range("E35").Value = '  +0.40%  '
$ws.Range("D36").Value = '0.01898'
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("E37").Value = '  +1.38%  '
$ws.Range("D38").Value = '0.8738'
$ws.Range("E38").Value = '  -4.60%  '
$ws.Range("D39").Value = '1.975'
$ws.Range("E39").Value = '  +0.20%  '
$ws.Range("D40").Value = '105.34'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '0.4073'
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("D43").Value = '5.504'
$ws.Range("E43").Value = '  -1.49%  '
$ws.Range("D44").Value = '7.207'
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("D45").Value = '62.17'
$ws.Range("E45").Value = '  +1.83%  '
$ws.Range("D46").Value = '0.1210'
$ws.Range("E46").Value = '  +0.97%  '
$ws.Range("D47").Value = '8.548'
$ws.Range("E47").Value = '  -1.74%  '
$ws.Range("D48").Value = '33.32'
$ws.Range("E48").Value = '  -0.64%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.05540'
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.373'
$ws.Range("E50").Value = '  -2.52%  '
$ws.Range("D51").Value = '0.3679'
$ws.Range("E51").Value = '  -0.74%  '
